$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = '29.903.45'
$ws.Range("E2").Value2 = '  +0.06%  '

# Row 3
$ws.Range("D3").Value2 = '1.901.73'
$ws.Range("E3").Value2 = '  +0.30%  '

# Row 4
$ws.Range("D4").Value2 = '''0.9996'
$ws.Range("E4").Value2 = '  -0.13%  '

# Row 5
$ws.Range("D5").Value2 = '''0.8022'
$ws.Range("E5").Value2 = '  +6.27%  '

# Row 6
$ws.Range("D6").Value2 = '''241.20'
$ws.Range("E6").Value2 = '  +0.42%  '

# Row 7
$ws.Range("D7").Value2 = '''0.9997'
$ws.Range("E7").Value2 = '  -0.11%  '

# Row 8
$ws.Range("D8").Value2 = '''0.3126'
$ws.Range("E8").Value2 = '  +2.87%  '

# Row 9
$ws.Range("D9").Value2 = '''26.18'
$ws.Range("E9").Value2 = '  +3.27%  '

# Row 10
$ws.Range("D10").Value2 = '''0.06882'

# Row 11
$ws.Range("D11").Value2 = '''0.07982'
$ws.Range("E11").Value2 = '  +0.16%  '

# Row 12
$ws.Range("D12").Value2 = '1.921.74'
$ws.Range("E12").Value2 = '  +1.18%  '

# Row 13
$ws.Range("E13").Value2 = '  -1.26%  '

# Row 14
$ws.Range("D14").Value2 = '''5.172'
$ws.Range("E14").Value2 = '  -0.48%  '

# Row 15
$ws.Range("D15").Value2 = '''92.42'
$ws.Range("E15").Value2 = '  +1.47%  '

# Row 16
$ws.Range("D16").Value2 = '29.900.05'
$ws.Range("E16").Value2 = '  +0.05%  '

# Row 17
$ws.Range("D17").Value2 = '''13.92'
$ws.Range("E17").Value2 = '  +0.10%  '

# Row 18
$ws.Range("E18").Value2 = '  -1.60%  '

# Row 19
$ws.Range("D19").Value2 = '''244.58'
$ws.Range("E19").Value2 = '  +0.66%  '

# Row 20
$ws.Range("D20").Value2 = '''0.000007696'
$ws.Range("E20").Value2 = '  -0.29%  '

# Row 21
$ws.Range("E21").Value2 = '  -0.05%  '

# Row 22
$ws.Range("D22").Value2 = '2.149.76'
$ws.Range("E22").Value2 = '  -0.46%  '

# Row 23
$ws.Range("D23").Value2 = '''0.9994'
$ws.Range("E23").Value2 = '  -0.19%  '

# Row 24
$ws.Range("D24").Value2 = '''6.897'
$ws.Range("E24").Value2 = '  -0.49%  '

# Row 25
$ws.Range("D25").Value2 = '''167.56'
$ws.Range("E25").Value2 = '  +1.27%  '

# Row 26
$ws.Range("D26").Value2 = '''9.180'
$ws.Range("E26").Value2 = '  -0.45%  '

# Row 27
$ws.Range("D27").Value2 = '''0.1425'
$ws.Range("E27").Value2 = '  +9.60%  '

# Row 28
$ws.Range("E28").Value2 = '  +0.54%  '

# Row 29
$ws.Range("D29").Value2 = '''2.033'
$ws.Range("E29").Value2 = '  +0.69%  '

# Row 30
$ws.Range("D30").Value2 = '''1.354'
$ws.Range("E30").Value2 = '  -3.80%  '

# Row 31
$ws.Range("D31").Value2 = '''1.512'
$ws.Range("E31").Value2 = '  -0.16%  '

# Row 32
$ws.Range("D32").Value2 = '''4.286'
$ws.Range("E32").Value2 = '  +0.34%  '

# Row 33
$ws.Range("D33").Value2 = '''0.05574'
$ws.Range("E33").Value2 = '  +4.12%  '

# Row 34
$ws.Range("E34").Value2 = '  +1.06%  '

# Row 35
$ws.Range("D35").Value2 = '''1.258'
$ws.Range("E35").Value2 = '  +0.96%  '

# Row 36
$ws.Range("D36").Value2 = '''0.7279'
$ws.Range("E36").Value2 = '  +0.56%  '

# Row 37
$ws.Range("D37").Value2 = '''2.720'
$ws.Range("E37").Value2 = '  +0.07%  '

# Row 38
$ws.Range("D38").Value2 = '''0.01931'
$ws.Range("E38").Value2 = '  +0.92%  '

# Row 39
$ws.Range("D39").Value2 = '''2.779'
$ws.Range("E39").Value2 = '  -0.21%  '

# Row 40
$ws.Range("D40").Value2 = '''0.4395'
$ws.Range("E40").Value2 = '  -0.04%  '

# Row 41
$ws.Range("B41").Value2 = 'FraxShare'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value2 = '''5.994'
$ws.Range("E41").Value2 = '  -2.82%  '

# Row 42
$ws.Range("B42").Value2 = 'Aave'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value2 = '''71.96'
$ws.Range("E42").Value2 = '  -0.26%  '

# Row 43
$ws.Range("D43").Value2 = '''0.9997'
$ws.Range("E43").Value2 = '  -0.16%  '

# Row 44
$ws.Range("D44").Value2 = '''0.8358'
$ws.Range("E44").Value2 = '  +1.36%  '

# Row 45
$ws.Range("D45").Value2 = '''1.860'
$ws.Range("E45").Value2 = '  -2.28%  '

# Row 46
$ws.Range("D46").Value2 = '''100.47'
$ws.Range("E46").Value2 = '  -0.63%  '

# Row 47
$ws.Range("D47").Value2 = '''7.554'
$ws.Range("E47").Value2 = '  +0.14%  '

# Row 48
$ws.Range("D48").Value2 = '''9.709'
$ws.Range("E48").Value2 = '  -0.44%  '

# Row 49
$ws.Range("D49").Value2 = '''977.00'
$ws.Range("E49").Value2 = '  +7.89%  '

# Row 50
$ws.Range("D50").Value2 = '2.058.48'
$ws.Range("E50").Value2 = '  +0.09%  '

# Row 51
$ws.Range("D51").Value2 = '''36.13'
$ws.Range("E51").Value2 = '  -0.08%  '
